# "Layout de sélection d'heure + travailler + dormir"
#
# Adds two new journal entries to the "Iteration #1" sheet's time-log table:
#   row 20: 2017-02-06 - "Création d'un layout pour chosir le nombre d'heure." - 3h
#   row 21: 2017-02-07 - "Peut travailler et dormir."                          - 2h
#
# The TOTAL cell (C36, =SUM(C14:C35)) recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration #1")

# Row 20
$ws.Cells.Item(20, 1).Value = 42772
$ws.Cells.Item(20, 2).Value = "Création d'un layout pour chosir le nombre d'heure."
$ws.Cells.Item(20, 3).Value = 3

# Row 21
$ws.Cells.Item(21, 1).Value = 42773
$ws.Cells.Item(21, 2).Value = "Peut travailler et dormir."
$ws.Cells.Item(21, 3).Value = 2

# These two rows were previously blank placeholder rows (date column styled
# for an empty cell). Now that they hold real dates, give column A the same
# look as the rows directly above it (row 18) by copying its format over.
$ws.Cells.Item(18, 1).Copy()
$ws.Cells.Item(20, 1).PasteSpecial(-4122)
$ws.Cells.Item(21, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0
